# "Generate Report for Handoff"
#
# b.md has now been handed off (in zh-cn and de-de) with a new handoff
# package (b.63290e5768f688058c7b37413b0a5c26c308f864.*.xlf), so its
# status flips from "Handed back: in sync with en-US" to
# "Ready for handoff", the handoff timestamp/file/datetime move forward,
# it is no longer flagged as a content duplicate, and an error is
# recorded noting the handback file used for a.md is stale relative to
# the newly generated b.md handoff.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet - row for b.md: zh-cn / de-de status + HO generate date
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2017-02-21 04:03:21"

# ---------------------------------------------------------------------
# zh-cn sheet - row for b.md
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("F3").Value = "'False"
$zhcn.Range("F3").Style = "Normal"
$zhcn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("H3").Value = "2017-02-21 04:03:06"
$zhcn.Range("R3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test4/blob/ba5d5f08b8e949aedf67f06ee9810d40859addf6/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test4/blob/5fa97ec0429b924ec8f0aa7602e88b45e84980cd/e2e/b.md."
$zhcn.Columns.Item(18).ColumnWidth = 39.17

# ---------------------------------------------------------------------
# de-de sheet - row for b.md
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("F3").Value = "'False"
$dede.Range("F3").Style = "Normal"
$dede.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("H3").Value = "2017-02-21 04:03:21"
$dede.Range("R3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test4/blob/ba5d5f08b8e949aedf67f06ee9810d40859addf6/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test4/blob/5fa97ec0429b924ec8f0aa7602e88b45e84980cd/e2e/b.md."
$dede.Columns.Item(18).ColumnWidth = 39.17
